# Applies the commit: inserts two new daily price rows (row 233 and 234)
# into the "Fruta, Terminal Hortofrutícola Agro Chillán - Naranja" sheet,
# shifting the existing rows 233-336 down to 235-338.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at the top of the data block (row 233), pushing
# all the existing rows (233-336) down by two (to 235-338).
$ws.Rows("233:234").Insert()

# --- New row 233: Naranja, Valencia, Primera, 2022-02-17 ---
$ws.Cells.Item(233, 1).Value = 7
$ws.Cells.Item(233, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(233, 3).Value = "Ñuble"
$ws.Cells.Item(233, 4).Value = 44609
$ws.Cells.Item(233, 5).Value = 16
$ws.Cells.Item(233, 6).Value = "Fruta"
$ws.Cells.Item(233, 7).Value = 100102
$ws.Cells.Item(233, 8).Value = "Cítricos"
$ws.Cells.Item(233, 9).Value = 100102005
$ws.Cells.Item(233, 10).Value = "Naranja"
$ws.Cells.Item(233, 11).Value = "Valencia"
$ws.Cells.Item(233, 12).Value = "Primera"
$ws.Cells.Item(233, 13).Value = 200
$ws.Cells.Item(233, 14).Value = 9500
$ws.Cells.Item(233, 15).Value = 10000
$ws.Cells.Item(233, 16).Value = 9750
$ws.Cells.Item(233, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(233, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(233, 19).Value = 650
$ws.Cells.Item(233, 20).Value = 15

# --- New row 234: Naranja, Valencia, Segunda, 2022-02-17 ---
$ws.Cells.Item(234, 1).Value = 7
$ws.Cells.Item(234, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(234, 3).Value = "Ñuble"
$ws.Cells.Item(234, 4).Value = 44609
$ws.Cells.Item(234, 5).Value = 16
$ws.Cells.Item(234, 6).Value = "Fruta"
$ws.Cells.Item(234, 7).Value = 100102
$ws.Cells.Item(234, 8).Value = "Cítricos"
$ws.Cells.Item(234, 9).Value = 100102005
$ws.Cells.Item(234, 10).Value = "Naranja"
$ws.Cells.Item(234, 11).Value = "Valencia"
$ws.Cells.Item(234, 12).Value = "Segunda"
$ws.Cells.Item(234, 13).Value = 60
$ws.Cells.Item(234, 14).Value = 8500
$ws.Cells.Item(234, 15).Value = 8500
$ws.Cells.Item(234, 16).Value = 8500
$ws.Cells.Item(234, 17).Value = "`$/bandeja 15 kilos granel"
$ws.Cells.Item(234, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(234, 19).Value = 567
$ws.Cells.Item(234, 20).Value = 15
